# Applies the "output generated at c8c62b6" update:
#  - bump the header date
#  - replace each multiplication problem's operands/result text
#
# Most of the problem cells have unique text in the document, so a simple
# Find/Replace (whole-document, match-case, no wildcards) is safe for them.
# The text "454×9=" appears twice (table row 5 col 3, and table row 20
# col 1) and maps to two different replacements, so those two cells are
# addressed directly via the Tables collection to avoid ambiguity.

$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Unique "2025-05-14 Wednesday" "2025-05-15 Thursday"

# Problem cells with text that is unique across the document
Replace-Unique "141×2=" "556×9="
Replace-Unique "756×9=" "667×3="
Replace-Unique "269×4=" "325×7="
Replace-Unique "287×5=" "475×3="
Replace-Unique "493×9=" "462×5="
Replace-Unique "915×6=" "197×9="
Replace-Unique "302×4=" "260×5="
Replace-Unique "602×9=" "681×3="
Replace-Unique "933×9=" "977×3="
Replace-Unique "401×9=" "564×2="
Replace-Unique "634×8=" "678×2="
Replace-Unique "941×4=" "495×6="
Replace-Unique "526×9=" "668×8="
Replace-Unique "877×5=" "479×7="
Replace-Unique "916×8=" "361×5="
Replace-Unique "979×6=" "207×2="
Replace-Unique "654×3=" "515×8="
Replace-Unique "223×3=" "186×9="
Replace-Unique "987×6=" "502×2="
Replace-Unique "830×8=" "198×4="
Replace-Unique "649×6=" "901×8="
Replace-Unique "177×9=" "792×5="
Replace-Unique "958×8=" "220×7="

# "454×9=" is duplicated (table row 5 col 3, and table row 20 col 1) with
# two different replacements, so target each cell explicitly.
$tbl = $d.Tables.Item(1)
$tbl.Cell(5, 3).Range.Text = "146×4="
$tbl.Cell(20, 1).Range.Text = "640×5="

Write-Host "Edit applied"
